# Updated cryptos list on Sat Dec 23 11:16:40 UTC 2023 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns of the crypto table with
# newly-scraped values, and re-ranks two coin pairs that swapped rank order
# (Celestia/MultiversX at rows 42-43, Aave/TrustWalletToken at rows 48-49).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D/E hold plain text in the source data (e.g. "43.588.44" uses
# dots as thousands separators, and some single-decimal prices like
# "94.98" would otherwise be auto-parsed as numbers by Excel). Force the
# Price cells that look numeric to Text format before writing so they are
# stored as literal strings, matching the original file's representation.

$ws.Range("D2").Value = "43.588.44"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "2.283.15"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "94.98"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.36"
$ws.Range("E6").Value = "  -1.45%  "
$ws.Range("E7").Value = "  -1.03%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.608"
$ws.Range("E9").Value = "  -2.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.34"
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0931"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.91"
$ws.Range("E12").Value = "  -3.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.105"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").Value = "2.626.22"
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.25"
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.845"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("D17").Value = "2.289.25"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("D18").Value = "43.561.16"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000110"
$ws.Range("E19").Value = "  +3.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.21"
$ws.Range("E20").Value = "  -1.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.94"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.61"
$ws.Range("E22").Value = "  +14.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.24"
$ws.Range("E23").Value = "  -2.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.07"
$ws.Range("E24").Value = "  -6.31%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.53"
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.17"
$ws.Range("E27").Value = "  -1.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.46"
$ws.Range("E28").Value = "  +2.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.01"
$ws.Range("E29").Value = "  +2.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.28"
$ws.Range("E30").Value = "  -4.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.86"
$ws.Range("E31").Value = "  +1.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.77"
$ws.Range("E32").Value = "  -3.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0894"
$ws.Range("E33").Value = "  -0.97%  "
$ws.Range("E34").Value = "  -4.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.125"
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.107"
$ws.Range("E36").Value = "  -3.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0352"
$ws.Range("E37").Value = "  -1.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.35"
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.34"
$ws.Range("E39").Value = "  -3.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.239"
$ws.Range("E40").Value = "  +1.73%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.72"
$ws.Range("E42").Value = "  +6.39%  "
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.27"
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("E44").Value = "  +1.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.75"
$ws.Range("E45").Value = "  -2.08%  "
$ws.Range("E46").Value = "  -5.15%  "
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.19"
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "96.20"
$ws.Range("E49").Value = "  -4.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.185"
$ws.Range("E50").Value = "  +8.38%  "
$ws.Range("D51").Value = "2.507.53"
$ws.Range("E51").Value = "  -0.88%  "
